$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E (value/volume columns) to be treated as text so that
# numeric-looking strings (e.g. "212.68", "1.00") are not auto-converted to numbers,
# matching the original workbook where these cells are inline/shared strings.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.550.40'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.639.17'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '212.68'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = '0.533'
$ws.Range('E6').Value = '  +4.90%  '
$ws.Range('D8').Value = '22.96'
$ws.Range('E8').Value = '  -5.33%  '
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').Value = '1.871.69'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = '1.639.34'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').Value = '0.563'
$ws.Range('E15').Value = '  -2.22%  '
$ws.Range('D16').Value = '64.33'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '27.525.18'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').Value = '230.24'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').Value = '7.71'
$ws.Range('E19').Value = '  +3.41%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '4.31'
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').Value = '9.91'
$ws.Range('E23').Value = '  +6.45%  '
$ws.Range('E24').Value = '  -3.71%  '
$ws.Range('D25').Value = '149.54'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('E27').Value = '  +1.58%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = '15.58'
$ws.Range('E29').Value = '  -2.90%  '
$ws.Range('D30').Value = '1.18'
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('E33').Value = '  +2.20%  '
$ws.Range('D34').Value = '1.423.09'
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('D35').Value = '1.58'
$ws.Range('E35').Value = '  +2.03%  '
$ws.Range('E36').Value = '  -1.81%  '
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '0.877'
$ws.Range('E38').Value = '  -3.65%  '
$ws.Range('E39').Value = '  -1.90%  '
$ws.Range('D40').Value = '0.886'
$ws.Range('E40').Value = '  +12.68%  '
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('E45').Value = '  +1.51%  '
$ws.Range('D46').Value = '64.91'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').Value = '1.781.08'
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('D48').Value = '1.67'
$ws.Range('E48').Value = '  -3.14%  '
$ws.Range('D49').Value = '86.07'
$ws.Range('E49').Value = '  -2.85%  '
$ws.Range('D50').Value = '0.0₆0104'
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').Value = '0.0990'
$ws.Range('E51').Value = '  -2.20%  '

# Restore the default cell style so no stray style index is left on the cells.
$dataRange.Style = "Normal"
